$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2 = @{ "E" = 3; "F" = 1; "G" = 2.17744; "H" = 6.53232; "I" = 0.6849679343335416; "J" = 0.6849679343335416; "K" = 3; "L" = 1; "M" = 0.8894906666666667; "N" = 2.668472; "O" = 0.1316647303664537; "P" = 0.1368307050811258; "Q" = 1.936812557226667; "R" = 17.43131301504; "S" = 0.0901861183836925; "T" = 0.09372464541282077 }
    3 = @{ "E" = 3; "F" = 1; "G" = 2.17744; "H" = 6.53232; "I" = 0.6849679343335416; "J" = 0.6849679343335416; "K" = 3; "L" = 1; "M" = 1.929018; "N" = 5.787053999999999; "O" = 0.2855382797818778; "P" = 0.2967416106155693; "Q" = 4.20032095392; "R" = 37.80288858528; "S" = 0.1955845656753457; "T" = 0.2032584880541546 }
    4 = @{ "E" = 3; "F" = 1; "G" = 2.17744; "H" = 6.53232; "I" = 0.6849679343335416; "J" = 0.6849679343335416; "K" = 3; "L" = 1; "M" = 1.949498666666667; "N" = 5.848496; "O" = 0.2885698815236896; "P" = 0.2998921597618951; "Q" = 4.244916376746668; "R" = 38.20424739072; "S" = 0.1976611156581565; "T" = 0.2054165131949298 }
    5 = @{ "E" = 3; "F" = 1; "G" = 2.17744; "H" = 6.53232; "I" = 0.6849679343335416; "J" = 0.6849679343335416; "K" = 3; "L" = 1; "M" = 1.222540333333333; "N" = 3.667621; "O" = 0.180963611404333; "P" = 0.1880638685361299; "Q" = 2.662008223413334; "R" = 23.95807401072; "S" = 0.1239542710931637; "T" = 0.1288177195539676 }
    6 = @{ "E" = 3; "F" = 1; "G" = 2.17744; "H" = 6.53232; "I" = 0.6849679343335416; "J" = 0.6849679343335416; "K" = 2; "L" = 1; "M" = 0.765177; "N" = 1.530354; "O" = 0.1132634969236461; "P" = 0.07847165600527987; "Q" = 1.66612700688; "R" = 9.99676204128; "S" = 0.07758186352318332; "T" = 0.05375056811766881 }
    7 = @{ "E" = 3; "F" = 1; "G" = 1.001453333333333; "H" = 3.00436; "I" = 0.3150320656664583; "J" = 0.3150320656664584; "K" = 3; "L" = 1; "M" = 0.8894906666666667; "N" = 2.668472; "O" = 0.1316647303664537; "P" = 0.1368307050811258; "Q" = 0.8907833931022222; "R" = 8.01705053792; "S" = 0.04147861198276116; "T" = 0.04310605966830502 }
    8 = @{ "E" = 3; "F" = 1; "G" = 1.001453333333333; "H" = 3.00436; "I" = 0.3150320656664583; "J" = 0.3150320656664584; "K" = 3; "L" = 1; "M" = 1.929018; "N" = 5.787053999999999; "O" = 0.2855382797818778; "P" = 0.2967416106155693; "Q" = 1.93182150616; "R" = 17.38639355544; "S" = 0.08995371410653208; "T" = 0.09348312256141462 }
    9 = @{ "E" = 3; "F" = 1; "G" = 1.001453333333333; "H" = 3.00436; "I" = 0.3150320656664583; "J" = 0.3150320656664584; "K" = 3; "L" = 1; "M" = 1.949498666666667; "N" = 5.848496; "O" = 0.2885698815236896; "P" = 0.2998921597618951; "Q" = 1.952331938062222; "R" = 17.57098744256; "S" = 0.09090876586553306; "T" = 0.09447564656696537 }
    10 = @{ "E" = 3; "F" = 1; "G" = 1.001453333333333; "H" = 3.00436; "I" = 0.3150320656664583; "J" = 0.3150320656664584; "K" = 3; "L" = 1; "M" = 1.222540333333333; "N" = 3.667621; "O" = 0.180963611404333; "P" = 0.1880638685361299; "Q" = 1.224317091951111; "R" = 11.01885382756; "S" = 0.05700934031116926; "T" = 0.05924614898216227 }
    11 = @{ "E" = 3; "F" = 1; "G" = 1.001453333333333; "H" = 3.00436; "I" = 0.3150320656664583; "J" = 0.3150320656664584; "K" = 2; "L" = 1; "M" = 0.765177; "N" = 1.530354; "O" = 0.1132634969236461; "P" = 0.07847165600527987; "Q" = 0.76628905724; "R" = 4.59773434344; "S" = 0.03568163340046278; "T" = 0.02472108788761106 }
}

foreach ($r in $rowData.Keys) {
    $cols = $rowData[$r]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$r").Value = $cols[$col]
    }
}
